$d = $word.ActiveDocument

# The edit relocates the hidden "_GoBack" bookmark from the end of the
# "...simultaneously." paragraph (Performance Requirements) down to a new
# empty paragraph that replaces the "If system crashes..." sentence
# (Security Requirements), and removes that sentence (and its redundant
# trailing blank paragraph) entirely.

# 1. Remove the existing _GoBack bookmark (it currently sits right after
#    "...simultaneously.").
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Find the paragraph containing the sentence to be removed and delete
#    the whole paragraph (text + paragraph mark). This merges it with the
#    following blank paragraph, leaving a single plain empty paragraph in
#    its place - matching the following blank paragraph that was already
#    there.
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -like "If system crashes it will return back at most one hour in maintainability purposes.*") {
        $para.Range.Delete()
        $found = $true
        break
    }
}

if (-not $found) {
    throw "Could not locate the 'If system crashes...' paragraph"
}

# 3. Re-add the _GoBack bookmark as a zero-length bookmark on the now
#    empty paragraph that used to hold that sentence.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    $prevText = ""
    if ($i -gt 1) {
        $prevText = $d.Paragraphs($i - 1).Range.Text
    }
    if ($para.Range.Text -eq [string][char]13 -and $prevText -like "*Passwords of the user shall be encrypted in DBMS for security purposes.*") {
        $d.Bookmarks.Add("_GoBack", $para.Range)
        break
    }
}
